$wb = $excel.ActiveWorkbook

# --- Sheet "List" (sheet1): drop the id and account columns ---
$ws1 = $wb.Worksheets.Item("List")

$ws1.Range("A1").Value = "`${msg.getProperty('savedSearch_name')}"
$ws1.Range("B1").Value = "`${msg.getProperty('savedSearch_formClassname')}"
$ws1.Range("C1").Value = "`${msg.getProperty('savedSearch_formContent')}"

$ws1.Range("A2").Value = "`${savedSearch.name}"
$ws1.Range("B2").Value = "`${savedSearch.formClassname}"
$ws1.Range("C2").Value = "`${savedSearch.formContent}"

$ws1.Range("D1:E2").ClearContents()

# --- Sheet "Search" (sheet2): shift key/value rows, rename account -> user, add id row ---
$ws2 = $wb.Worksheets.Item("Search")

$ws2.Range("A4").Value = "`${msg.getProperty('savedSearch_id')}"
$ws2.Range("B4").Value = "`${id}"

$ws2.Range("A5").Value = "`${msg.getProperty('savedSearch_name')}"
$ws2.Range("B5").Value = "`${name}"

$ws2.Range("A6").Value = "`${msg.getProperty('savedSearch_formClassname')}"
$ws2.Range("B6").Value = "`${formClassname}"

$ws2.Range("A7").Value = "`${msg.getProperty('savedSearch_formContent')}"
$ws2.Range("B7").Value = "`${formContent}"

$ws2.Range("A8").Value = "`${msg.getProperty('savedSearch_user')}"
$ws2.Range("B8").Value = "`${user}"
